$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.71"
$ws.Range("E2").Value = "'-6.94%"
$ws.Range("D3").Value = "'40.74"
$ws.Range("E3").Value = "'-0.30%"
$ws.Range("D4").Value = "'5.038"
$ws.Range("E4").Value = "'-1.76%"
$ws.Range("D5").Value = "'0.07399"
$ws.Range("E5").Value = "'-3.03%"
$ws.Range("D6").Value = "'4.278"
$ws.Range("E6").Value = "'-1.34%"
$ws.Range("D7").Value = "'1.553"
$ws.Range("E7").Value = "'-7.82%"
$ws.Range("D8").Value = "'0.9252"
$ws.Range("E8").Value = "'-1.02%"
$ws.Range("D9").Value = "'0.1155"
$ws.Range("E9").Value = "'-7.61%"
$ws.Range("D10").Value = "'0.1727"
$ws.Range("E10").Value = "'-5.03%"
$ws.Range("D11").Value = "'0.08664"
$ws.Range("E11").Value = "'-3.96%"
$ws.Range("D12").Value = "'0.04171"
$ws.Range("E12").Value = "'0.73%"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("D14").Value = "'0.001273"
$ws.Range("E14").Value = "'-1.65%"
$ws.Range("D15").Value = "'0.005918"
$ws.Range("E15").Value = "'1.88%"
$ws.Range("D16").Value = "'3.413"
$ws.Range("E16").Value = "'1.74%"
$ws.Range("E18").Value = "'-2.26%"
$ws.Range("D19").Value = "'7.641"
$ws.Range("E19").Value = "'-9.20%"
$ws.Range("D20").Value = "'0.1377"
$ws.Range("E20").Value = "'2.08%"
$ws.Range("D21").Value = "'0.2875"
$ws.Range("E21").Value = "'4.86%"
$ws.Range("D22").Value = "'0.03856"
$ws.Range("E22").Value = "'-4.54%"
$ws.Range("E23").Value = "'-0.74%"
$ws.Range("D24").Value = "'0.003856"
$ws.Range("E24").Value = "'-4.85%"
$ws.Range("D25").Value = "'0.0001277"
$ws.Range("E25").Value = "'0.15%"
$ws.Range("D26").Value = "'0.0003714"
$ws.Range("D38").Value = "'0.02344"
$ws.Range("E38").Value = "'-5.58%"
$ws.Range("D39").Value = "'0.05018"
$ws.Range("E39").Value = "'-3.37%"
$ws.Range("D40").Value = "'0.005853"
$ws.Range("E40").Value = "'170.73%"
$ws.Range("D41").Value = "'0.007670"
$ws.Range("E41").Value = "'-1.03%"
$ws.Range("D42").Value = "'0.1286"
$ws.Range("E42").Value = "'-1.07%"
$ws.Range("D43").Value = "'0.007329"
$ws.Range("E43").Value = "'-0.68%"
$ws.Range("D44").Value = "'0.007098"
$ws.Range("E44").Value = "'-13.67%"
$ws.Range("D45").Value = "'0.3158"
$ws.Range("E45").Value = "'0.65%"
$ws.Range("D46").Value = "'0.00006408"
$ws.Range("E46").Value = "'-3.94%"
$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.75%"
$ws.Range("E48").Value = "'-93.30%"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("E49").Value = "'-0.75%"
$ws.Range("D50").Value = "'0.0001996"
$ws.Range("E50").Value = "'-0.75%"
